{"js": "// 1. Update the street address: \"#401\" -> \"Suite 401\"\nconst addrResults = context.document.body.search(\"133 Hazelton Ave, #401\", { matchCase: true, matchWholeWord: false });\naddrResults.load(\"text\");\nawait context.sync();\nif (addrResults.items.length > 0) {\n  addrResults.items[0].insertText(\"133 Hazelton Ave, Suite 401\", Word.InsertLocation.replace);\n}\n\n// 2. Underline the email address hyperlink run\nconst emailResults = context.document.body.search(\"mikecorbridge@gmail.com\", { matchCase: true, matchWholeWord: false });\nemailResults.load(\"font\");\nawait context.sync();\nfor (let i = 0; i < emailResults.items.length; i++) {\n  emailResults.items[i].font.underline = Word.UnderlineType.single;\n}\n\n// 3. Re-color three project/employer hyperlink URLs from black to blue\nconst blueLinks = [\n  \"http://www.twintechs.com/\",\n  \"http://www.sound-shapes.com/\",\n  \"http://nexthabitatadvisors.com/\"\n];\nfor (const url of blueLinks) {\n  const linkResults = context.document.body.search(url, { matchCase: true, matchWholeWord: false });\n  linkResults.load(\"font\");\n  await context.sync();\n  for (let i = 0; i < linkResults.items.length; i++) {\n    linkResults.items[i].font.color = \"#0000FF\";\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the street address: \"#401\" -> \"Suite 401\"\n$rng = $d.Content\n$rng.Find.Text = \"133 Hazelton Ave, #401\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\nif ($rng.Find.Execute()) {\n    # Replace on the owning paragraph's Range (rather than the Find range\n    # itself) so the run's xml:space=\"preserve\" text node is edited in\n    # place instead of being rebuilt without it.\n    $pr = $rng.Paragraphs(1).Range\n    $pr.Text = \"133 Hazelton Ave, Suite 401\"\n}\n\n# 2. Underline the email address hyperlink run\n$rng = $d.Content\n$rng.Find.Text = \"mikecorbridge@gmail.com\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\nif ($rng.Find.Execute()) {\n    $rng.Font.Underline = 1  # wdUnderlineSingle\n}\n\n# 3. Re-color three project/employer hyperlink URLs from black to blue\n$blueLinks = @(\n    \"http://www.twintechs.com/\",\n    \"http://www.sound-shapes.com/\",\n    \"http://nexthabitatadvisors.com/\"\n)\nforeach ($url in $blueLinks) {\n    $rng = $d.Content\n    $rng.Find.Text = $url\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    if ($rng.Find.Execute()) {\n        $rng.Font.Color = \"0000FF\"\n    }\n}\n"}
